$d = $word.ActiveDocument

# --- Paragraph 1: title "Descrição das abas do excel." ---------------------
# Remove the spell-check proofErr markers and capitalize "excel" -> "Excel"
# by deleting the exact span (which fully covers the proofErr-wrapped word)
# and retyping it.
$titleFix = $d.Content
$titleFix.Find.Execute("do excel.", $true, $false, $false, $false, $false, `
    $true, 1, $false, "", 0) | Out-Null
$titleFix.Delete()
$titleFix.InsertAfter("do Excel.")

# Bold + 20pt the title, applied to three contiguous sub-ranges ("Descrição
# das abas do ", "E", "xcel") so each keeps its own run, then leave the
# trailing "." unformatted.
$prefixRng = $d.Content
$prefixRng.Find.Execute("Descrição das abas do ", $true, $false, $false, `
    $false, $false, $true, 1, $false, "", 0) | Out-Null
$prefixRng.Font.Bold = $true
$prefixRng.Font.Size = 20

$eRng = $d.Range($prefixRng.End, $prefixRng.End + 1)
$eRng.Font.Bold = $true
$eRng.Font.Size = 20

$xcelRng = $d.Range($eRng.End, $eRng.End + 4)
$xcelRng.Font.Bold = $true
$xcelRng.Font.Size = 20

# --- Section labels: bold + underline the label, keep the rest plain -------
$labels = @(
    "Introdução",
    "Colaborar",
    "Funções e formas",
    "Importar e analisar",
    "Formatar dados",
    "Solução de problemas"
)

foreach ($label in $labels) {
    $rng = $d.Content
    $found = $rng.Find.Execute($label, $true, $false, $false, $false, $false, `
        $true, 1, $false, "", 0)
    if ($found) {
        $rng.Font.Bold = $true
        $rng.Font.Underline = 1
    }
}
